$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the skills table (rows 2-13) with the new "About me" data.
# Columns: A=id, B=topic, C=skills, D=level, E=tooltip

# Row 2 is special-cased: the tooltip (E2) was entered before the skill (C2)
# while editing, so set it first to reproduce the original shared-string order.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Programming languages"
$ws.Cells.Item(2, 5).Value = "1+ years experience"
$ws.Cells.Item(2, 3).Value = "Python"
$ws.Cells.Item(2, 4).Value = 4

# row, id, topic, skill, level, tooltip
$data = @(
    @(3,  2,  "Programming languages", "Java",               2, "<1 year experience"),
    @(4,  3,  "Programming languages", "R",                   1, "<1 year experience"),
    @(5,  4,  "Machine learning",      "Supervised ML",       4, "Classification & regression, time series"),
    @(6,  5,  "Machine learning",      "Unsupervised ML",     3, "K-Means, PCA"),
    @(7,  6,  "Machine learning",      "Deep Learning",       3, "TensorFlow"),
    @(8,  7,  "Software engineering",  "Front-end",           2, "HTML, CSS"),
    @(9,  8,  "Software engineering",  "Databases",           3, "Relational and NoSQL"),
    @(10, 9,  "Software engineering",  "Deployment",          3, "Heroku, AWS"),
    @(11, 10, "Data analysis",         "Data Wrangling",      4, "Pandas, Numpy"),
    @(12, 11, "Data analysis",         "Data Visuzlization",  4, "Matplotlib, Seaborn, Plotly, Shapley"),
    @(13, 12, "Data analysis",         "Statistics",          3, "Descriptive, inferential")
)

foreach ($entry in $data) {
    $row = $entry[0]
    $id = $entry[1]
    $topic = $entry[2]
    $skill = $entry[3]
    $level = $entry[4]
    $tooltip = $entry[5]

    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $topic
    $ws.Cells.Item($row, 3).Value = $skill
    $ws.Cells.Item($row, 4).Value = $level
    $ws.Cells.Item($row, 5).Value = $tooltip
}

$ws.Range("E14").Select()
